$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 1.85
$ws.Range("H5").Value = 3.7
$ws.Range("I5").Value = 3.9
$ws.Range("J5").Value = 2.6
$ws.Range("S5").Value = 1.5
$ws.Range("T5").Value = 2.5
$ws.Range("X5").Value = 8
$ws.Range("AG5").Value = 9
$ws.Range("AH5").Value = 19
$ws.Range("AI5").Value = 13
$ws.Range("AK5").Value = 34
$ws.Range("AR5").Value = 67
$ws.Range("AT5").Value = 2.5
$ws.Range("AZ5").Value = 81
$ws.Range("J6").Value = 2.75
$ws.Range("L6").Value = 3.6
$ws.Range("N6").Value = 7.9
$ws.Range("O6").Value = 1.31
$ws.Range("P6").Value = 2.9
$ws.Range("Q6").Value = 1.95
$ws.Range("R6").Value = 1.75
$ws.Range("V6").Value = 1.91
$ws.Range("W6").Value = 7.6
$ws.Range("X6").Value = 11
$ws.Range("Z6").Value = 23
$ws.Range("AA6").Value = 18
$ws.Range("AC6").Value = 9
$ws.Range("AD6").Value = 6.1
$ws.Range("AG6").Value = 9.25
$ws.Range("AH6").Value = 16.5
$ws.Range("AI6").Value = 11
$ws.Range("AK6").Value = 28
$ws.Range("AL6").Value = 35
$ws.Range("AM6").Value = 500
$ws.Range("AP6").Value = 18
$ws.Range("AR6").Value = 70
$ws.Range("AT6").Value = 2.6
$ws.Range("AX6").Value = 16.5
$ws.Range("AY6").Value = 22
$ws.Range("AZ6").Value = 80
$ws.Range("BA6").Value = 110
$ws.Range("H7").Value = 4.6
$ws.Range("I7").Value = 7
$ws.Range("O7").Value = 1.14
$ws.Range("Q7").Value = 1.55
$ws.Range("R7").Value = 2.15
$ws.Range("W7").Value = 8
$ws.Range("X7").Value = 7.1
$ws.Range("Z7").Value = 9.25
$ws.Range("AE7").Value = 18.5
$ws.Range("AG7").Value = 21
$ws.Range("AL7").Value = 65
$ws.Range("AM7").Value = 600
$ws.Range("AP7").Value = 15
$ws.Range("AT7").Value = 3.2
$ws.Range("AU7").Value = 7.9
$ws.Range("AX7").Value = 40
$ws.Range("AY7").Value = 37
$ws.Range("BB7").Value = 450
$ws.Range("O11").Value = 1.3
$ws.Range("P11").Value = 3.4
$ws.Range("Q11").Value = 1.98
$ws.Range("R11").Value = 1.83
$ws.Range("G19").Value = 1.9
$ws.Range("I19").Value = 4.1
$ws.Range("J19").Value = 2.5
$ws.Range("L19").Value = 4.5
$ws.Range("W19").Value = 7
$ws.Range("Z19").Value = 15
$ws.Range("AG19").Value = 11
$ws.Range("AW19").Value = 6
$ws.Range("BA19").Value = 101
$ws.Range("G21").Value = 3.4
$ws.Range("I21").Value = 2.3
$ws.Range("J21").Value = 4
$ws.Range("L21").Value = 3.1
$ws.Range("AY21").Value = 29
$ws.Range("AZ21").Value = 51
$ws.Range("G23").Value = 2.6
$ws.Range("I23").Value = 2.45
$ws.Range("L23").Value = 3
$ws.Range("X23").Value = 15
$ws.Range("Z23").Value = 26
$ws.Range("AI23").Value = 10
$ws.Range("AJ23").Value = 26
